$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the grade mark "ок" in newly filled cells
$ws.Range("C5").Value = "ок"

$ws.Range("C12").Value = "ок"
$ws.Range("D12").Value = "ок"
$ws.Range("E12").Value = "ок"

$ws.Range("C16").Value = "ок"
$ws.Range("D16").Value = "ок"

$ws.Range("C20").Value = "ок"
$ws.Range("D20").Value = "ок"
$ws.Range("E20").Value = "ок"
$ws.Range("F20").Value = "ок"
$ws.Range("G20").Value = "ок"
$ws.Range("H20").Value = "ок"

# I20 is a brand-new cell outside the pre-existing table body, so give it
# the same formatting (thick border, centered, wrapped) as its neighbour
# H20 before writing the value into it.
$ws.Range("H20").Copy()
$ws.Range("I20").PasteSpecial(-4122)
$ws.Range("I20").Value = "ок"

$ws.Range("C29").Value = "ок"

# Update the view: scroll position and active cell selection
$window = $excel.ActiveWindow
$window.ScrollColumn = 3
$window.ScrollRow = 5
$window.RangeSelection.Worksheet.Range("I20").Select()
